# Insert a new row at position 73 (shifts existing rows 73:141 down to 74:142)
# and populate it with the new weekly price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("73:73").Insert()

$ws.Cells.Item(73, 1).Value  = 5
$ws.Cells.Item(73, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(73, 3).Value  = "Maule"
$ws.Cells.Item(73, 4).Value  = 45167
$ws.Cells.Item(73, 5).Value  = 7
$ws.Cells.Item(73, 6).Value  = 100112013
$ws.Cells.Item(73, 7).Value  = "Alcachofa"
$ws.Cells.Item(73, 8).Value  = "Madrigal"
$ws.Cells.Item(73, 9).Value  = "Primera"
$ws.Cells.Item(73, 10).Value = 300
$ws.Cells.Item(73, 11).Value = 13000
$ws.Cells.Item(73, 12).Value = 13000
$ws.Cells.Item(73, 13).Value = 13000
$ws.Cells.Item(73, 14).Value = '$/caja 40 unidades'
$ws.Cells.Item(73, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(73, 16).Value = 325
$ws.Cells.Item(73, 17).Value = 40
$ws.Cells.Item(73, 18).Value = "Hortaliza"
